$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-18 Monday" "2024-11-19 Tuesday"

Replace-Text "49÷2=" "65÷7="
Replace-Text "58÷4=" "90÷4="
Replace-Text "62÷3=" "55÷7="
Replace-Text "77÷3=" "37÷4="
Replace-Text "28÷6=" "20÷8="

Replace-Text "14÷3=" "93÷6="
Replace-Text "63÷6=" "23÷9="
Replace-Text "23÷2=" "76÷7="
Replace-Text "54÷4=" "51÷4="
Replace-Text "64÷9=" "63÷7="

Replace-Text "15÷5=" "53÷7="
Replace-Text "15÷4=" "34÷6="
Replace-Text "28÷4=" "96÷9="
Replace-Text "28÷2=" "22÷8="
Replace-Text "19÷5=" "95÷8="

Replace-Text "30÷2=" "78÷7="
Replace-Text "83÷9=" "93÷6="
Replace-Text "51÷5=" "52÷5="
Replace-Text "76÷2=" "66÷2="
Replace-Text "29÷4=" "39÷2="

Replace-Text "50÷2=" "71÷2="
Replace-Text "69÷3=" "66÷2="
Replace-Text "38÷7=" "78÷7="
Replace-Text "39÷4=" "99÷4="
Replace-Text "40÷2=" "81÷4="

Write-Output "done"
